$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.620.91"
$ws.Range("E2").Value = "  +1.29%  "

# Row 3
$ws.Range("D3").Value = "2.454.42"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.26%  "

# Row 5
$ws.Range("D5").Value = "'491.46"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6
$ws.Range("D6").Value = "'155.91"
$ws.Range("E6").Value = "  +2.14%  "

# Row 7
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").Value = "'0.612"
$ws.Range("E8").Value = "  +18.75%  "

# Row 9
$ws.Range("D9").Value = "2.492.77"
$ws.Range("E9").Value = "  +1.30%  "

# Row 10
$ws.Range("D10").Value = "'6.20"
$ws.Range("E10").Value = "  +8.73%  "

# Row 11
$ws.Range("D11").Value = "'0.101"
$ws.Range("E11").Value = "  +1.04%  "

# Row 12
$ws.Range("D12").Value = "'0.335"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
$ws.Range("E13").Value = "  +0.88%  "

# Row 14
$ws.Range("D14").Value = "2.882.27"
$ws.Range("E14").Value = "  +0.02%  "

# Row 15
$ws.Range("D15").Value = "57.531.84"
$ws.Range("E15").Value = "  +0.74%  "

# Row 16
$ws.Range("D16").Value = "'20.83"
$ws.Range("E16").Value = "  -0.75%  "

# Row 17
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = "  -2.25%  "

# Row 18
$ws.Range("D18").Value = "2.482.51"
$ws.Range("E18").Value = "  +0.54%  "

# Row 19
$ws.Range("D19").Value = "'4.65"
$ws.Range("E19").Value = "  +1.06%  "

# Row 20
$ws.Range("D20").Value = "'325.74"
$ws.Range("E20").Value = "  +0.24%  "

# Row 21
$ws.Range("D21").Value = "'10.11"
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("D22").Value = "'0.996"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").Value = "'5.96"
$ws.Range("E23").Value = "  +2.09%  "

# Row 24
$ws.Range("D24").Value = "'58.53"
$ws.Range("E24").Value = "  +0.70%  "

# Row 25
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.404"
$ws.Range("E25").Value = "  -0.70%  "

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.992"
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -0.79%  "

# Row 28
$ws.Range("D28").Value = "2.562.51"
$ws.Range("E28").Value = "  -0.20%  "

# Row 29
$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = "  -2.80%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0810"
$ws.Range("E30").Value = "  +0.51%  "

# Row 31
$ws.Range("E31").Value = "  -0.05%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'18.82"
$ws.Range("E32").Value = "  +3.27%  "

# Row 33
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'150.84"
$ws.Range("E33").Value = "  -0.14%  "

# Row 34
$ws.Range("D34").Value = "'1.53"
$ws.Range("E34").Value = "  +0.78%  "

# Row 35
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  +1.91%  "

# Row 36
$ws.Range("D36").Value = "'3.81"
$ws.Range("E36").Value = "  +1.40%  "

# Row 37
$ws.Range("D37").Value = "'1.14"
$ws.Range("E37").Value = "  -1.10%  "

# Row 38
$ws.Range("D38").Value = "'0.835"
$ws.Range("E38").Value = "  -6.18%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.40"
$ws.Range("E39").Value = "  -0.24%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'34.27"
$ws.Range("E40").Value = "  +0.16%  "

# Row 41
$ws.Range("D41").Value = "'3.55"
$ws.Range("E41").Value = "  +1.22%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'284.21"
$ws.Range("E42").Value = "  +6.80%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.100"
$ws.Range("E43").Value = "  +4.40%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.610"
$ws.Range("E44").Value = "  +0.40%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.994"
$ws.Range("E45").Value = "  -0.21%  "

# Row 46
$ws.Range("D46").Value = "'0.0538"
$ws.Range("E46").Value = "  -3.71%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0230"
$ws.Range("E47").Value = "  +0.31%  "

# Row 48
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.24"
$ws.Range("E48").Value = "  +0.17%  "

# Row 49
$ws.Range("D49").Value = "'4.67"
$ws.Range("E49").Value = "  -3.05%  "

# Row 50
$ws.Range("D50").Value = "'18.08"
$ws.Range("E50").Value = "  +1.51%  "

# Row 51
$ws.Range("D51").Value = "1.906.12"
$ws.Range("E51").Value = "  +4.09%  "
